$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/living-location-pre-stroke-extension"
$wsMeta.Range("B8").Value = "2023-08-16T00:27:03-03:00"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z5").Value = "https://molic-avc.gabriellesantosleandro.com/ValueSet/living-location-pre-stroke-valueset"
$wsElem.Columns.Item(26).ColumnWidth = 81.3
